# Weekly update: insert a new price record for "Camote" at row 5,
# pushing the existing rows 5-22 down to 6-23 (newest week's data goes
# right after the header block's 3rd entry, matching the source feed's
# insertion order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record; Excel shifts rows 5:22 down to 6:23
# and carries their formatting (incl. the date-style cell in column D).
$ws.Rows(5).Insert()

# Populate the newly inserted row 5 with this week's observation.
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 45222
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100114002
$ws.Range("G5").Value = "Camote"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16500
$ws.Range("N5").Value = "$/malla 18 kilos"
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 917
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"
